# Applies the "AB-Tests.xlsx" edit:
#  - rename table column "A (Sekunden)" -> "c" (header cell + table column)
#  - repoint the AVERAGE formula in A15 at the renamed column
#  - select A1:B15 on the "Reigster" sheet
#  - set the page setup (paper size / orientation) for the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reigster")

# Rename the first table column header from "A (Sekunden)" to "c".
# Writing the header cell's value keeps the ListObject / table definition
# (xl/tables/table1.xml) and the shared-strings table in sync.
$ws.Range("A1").Value = "c"

# The formula text isn't rewritten automatically when the referenced
# column is renamed, so update it explicitly to keep the structured
# reference valid.
$ws.Range("A15").Formula = "=AVERAGE(Tabelle1[c])"

# Update the view selection to span the whole table (A1:B15).
$ws.Range("A1:B15").Select()

# Configure the page setup (paper size = A4/"9", orientation = portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
